# TC03_PDP_Check.xlsx -- "Update with new data"
#
# The Testdata sheet's Baseurl test-data cell (B2) held a hyperlink to a
# Nike PDP URL. Replace it with a hyperlink to the Under Armour PDP URL,
# and update the active-sheet/selection UI state to match the saved file
# (Testdata tab active, with new selections on both sheets).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TC03_PDP_Check
$ws2 = $wb.Worksheets.Item(2)   # Testdata

$newUrl = "http://129.213.54.196:8002/under-armour-heatgear-sonic-fitted-shirt/6225774?"

# Replace the hyperlink + its displayed text on Testdata!B2.
$urlCell = $ws2.Range("B2")
$urlCell.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($urlCell, $newUrl, "", "", $newUrl)

# Match the saved selections / active sheet.
$ws1.Range("B10").Select()
$ws2.Range("G12").Select()
$ws2.Activate()
